$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing "general" text formatting for the Price (D) column
# cells whose new value happens to look like a pure number (e.g. "305.50",
# "1.00", "0.0801"). Excel auto-converts a numeric-looking string assigned
# via .Value into an actual number (dropping trailing/insignificant zeros),
# so those specific cells are forced to Text format first to keep the exact
# original string such as "305.50" or "1.00".

$ws.Range("D2").Value = "44.616.43"
$ws.Range("E2").Value = "  +1.63%  "

$ws.Range("D3").Value = "2.233.25"
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.50"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.30"
$ws.Range("E6").Value = "  -1.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -0.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  -2.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.73"
$ws.Range("E10").Value = "  -1.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.14"
$ws.Range("E12").Value = "  -1.47%  "

$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").Value = "2.575.62"
$ws.Range("E14").Value = "  -0.57%  "

$ws.Range("D15").Value = "2.233.26"
$ws.Range("E15").Value = "  -3.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.829"
$ws.Range("E16").Value = "  -0.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.49"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").Value = "44.376.61"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("D19").Value = "0.0₃0941"
$ws.Range("E19").Value = "  -3.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.82"
$ws.Range("E20").Value = "  -2.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.18"
$ws.Range("E21").Value = "  -3.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.92"
$ws.Range("E22").Value = "  -0.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.12"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  -1.12%  "

$ws.Range("E25").Value = "  -1.66%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.29"
$ws.Range("E27").Value = "  +3.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.70"
$ws.Range("E28").Value = "  -3.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.17"
$ws.Range("E29").Value = "  -1.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.87"
$ws.Range("E30").Value = "  -2.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.78"
$ws.Range("E31").Value = "  -1.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.05"
$ws.Range("E32").Value = "  -1.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0782"
$ws.Range("E33").Value = "  -2.54%  "

$ws.Range("E34").Value = "  +1.00%  "

$ws.Range("E35").Value = "  -6.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.107"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("E37").Value = "  -2.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.86"
$ws.Range("E38").Value = "  +5.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.13"
$ws.Range("E39").Value = "  +3.11%  "

$ws.Range("E40").Value = "  -4.06%  "

$ws.Range("E41").Value = "  -3.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0299"
$ws.Range("E42").Value = "  +0.34%  "

$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").Value = "1.822.05"
$ws.Range("E44").Value = "  +4.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.73"
$ws.Range("E45").Value = "  +10.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "79.23"
$ws.Range("E46").Value = "  -4.64%  "

$ws.Range("E47").Value = "  -2.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.94"
$ws.Range("E48").Value = "  -2.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.83"
$ws.Range("E49").Value = "  -2.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.71"
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.96"
$ws.Range("E51").Value = "  -2.56%  "
